$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# Insert a new header row at the top, pushing existing data down one row
$ws.Rows.Item(1).Insert()

# Populate the new header row with the same column headers used on "Sheet1"
$ws.Cells.Item(1, 1).Value = "Company"
$ws.Cells.Item(1, 2).Value = "Revenue"
$ws.Cells.Item(1, 3).Value = "Founding Year"
